{"js": "// Apply the report-refresh edit: updated timestamp plus five rewritten\n// section bodies (Company Overview, Recent News, Competitors, Insights).\n// `\\u000b` encodes a Word manual line break (<w:br/>); paragraphs are\n// replaced wholesale (including their internal line breaks) via\n// `insertText(..., \"Replace\")`, which preserves the paragraph itself\n// (and its paragraph style) while swapping its run content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  [2, \"Generated: 2025-08-06 14:02:06\"],\n  [5, \"As an AI developed by OpenAI, I can deliver information based on the data I was last trained on, which includes information up to September 2021. For this reason, I'm unable to provide real-time or precise company-specific information or updates that have occurred post-September 2021. However, I can provide a general response based on the details you have given above:\\u000b\\u000b1. JP Morgan Chase, one of the world's most valuable banks, provides services to millions of consumers, small businesses, and many of the world's most prominent corporate, institutional, and government clients. They're major players in investment banking, financial services for consumers and businesses, financial transaction processing, asset management, and private equity.\\u000b\\u000b2. The current geopolitical environment, marked by the global economic recovery from the Covid-19 pandemic, trade tensions, and regulatory changes, can impact the bank's business operations, risk profile, and market dynamics. These impacts can vary greatly depending on the specific circumstances.\\u000b\\u000b3. Amidst its competitors like Bank of America, Goldman Sachs, and Citigroup, JP Morgan Chase often stands out with its vast portfolio of services, robust financial performance, and a strong brand name. \\u000b\\u000b4. Strengths of JP Morgan include its strong market position, diversified business model, and robust capital position. It also has a strong technology infrastructure that offers online and mobile banking options to customers.\\u000b\\u000b5. Key challenges for JP Morgan Chase can include regulatory changes, geopolitical uncertainties, and economic uncertainties related to global events like the Covid-19 pandemic. They undertake comprehensive risk management strategies to mitigate these risks, but the effectiveness of these strategies can vary depending on the scope of external events.\\u000b\\u000b6. JP Morgan Chase distinguishes itself from competitors with its strong emphasis on digital banking and technologies, a diversified business model, and strong risk management capabilities.\\u000b\\u000b7. I don't have up-to-date data on recent projects or initiatives started by JP Morgan Chase. However, the company has been known for investing heavily in technological advancements and customer-centric initiatives.\\u000b\\u000b8. I am unable to provide specific forecast figures for future growth. It's generally recommended to consult specialized financial sources or professionals for accurate figures.\\u000b\\u000b9. Comparison with competitors on revenue growth forecasts would also require up-to-date data from financial sources or professionals.\\u000b\\u000b10. Currently, I can't provide highlights from the last investor day presentation as my last training data covers information up until September 2021.\\u000b\\u000b11. For detailed comparisons of the latest 10Qs, you should consult professional financial sources or regulatory disclosures. Even more so for the comparison with the latest 10Qs of other competitors.\\u000b\\u000b12. A digital transformation plan typically includes a broad range of initiatives that use digital technologies to modify or create new business processes, culture, and customer experiences to meet changing business and market requirements. As of my last training data, JP Morgan Chase has been actively investing in digital banking, but specific initiatives for the upcoming year are not known.\\u000b\\u000b13. Over the past years, JP Morgan Chase has implemented several digital transformation and operational efficiency initiatives. They have been investing heavily in technologies like AI and blockchain, improving their online and mobile banking platforms and customer-centric digital tools.\\u000b\\u000b14. The exact individuals responsible for digital transformation in JP Morgan Chase can change over time, but large corporations typically have a Chief Digital Officer and a team dedicated to digital transformation initiatives. They work in conjunction with other departments to ensure business-wide digital strategies align.\\u000b\\u000bI encourage you to consult official JP Morgan Chase resources or related regulatory filings for more updated and accurate data.\"],\n  [7, \"1. **Trump Claims Rejection by JP Morgan:** Donald Trump claimed that JP Morgan Chase and Bank of America rejected him as a customer, alleging discrimination (Source: CNBC, The New York Times, Politico, Financial Times). This could draw regulatory attention and potential legal implications for the firm, impacting its reputation and investor relations.\\u000b\\u000b2. **Employee Complaints about Fitness Center Costs:** JP Morgan staff have expressed dissatisfaction about having to pay for the new high-end fitness center at the company's headquarters, as they return to office fulltime (Source: New York Post). This could possibly impact employee morale and turnover rates and might require the company to revisit its employee engagement strategies.\\u000b\\u000b3. **Expansion in Maryland:** JP Morgan plans to open four new branches in Maryland, signaling its intent to expand its footprint in the region through 2025 (Source: Maryland Daily Record). This is a positive development and could contribute to overall business growth and market share increase.\\u000b\\u000b4. **Acquisition of 3M Co. Shares:** The company has significantly expanded its stake in the 3M Co. with the purchase of over 13 million shares (Source: GuruFocus). This investment could enhance JP Morgan's financial returns depending on the performance of 3M's stock.\\u000b\\u000b5. **Potential Apple Play:** There are speculations that JP Morgan might make a play for Apple's credit card portfolio, possibly signaling a greater ambition to conquer more business in technology-based financial services and products (Source: Payments Dive). This could symbolise strategic change, moving towards a more tech-oriented business model.\\u000b\"],\n  [11, \"1. Bank of America Corp:\\u000bBank of America is one of the world's largest financial institutions, serving individual consumers, small and middle-market businesses and large corporations with a broad spectrum of banking, investing, asset management and other financial services. Its finance and leadership segment is directed by CEO Brian Moynihan. One key differentiation between JP Morgan Chase and Bank of America is that the latter has a more diversified business mix and is heavily inclined towards consumer banking.\\u000b\\u000b2. Wells Fargo & Co:\\u000bWells Fargo is a diversified, community-based financial services company that provides banking, investment and mortgage products and services, as well as consumer and commercial finance. The company is helmed by CEO Charles Scharf. One key differentiator is its network of branches across the United States which is larger than JP Morgan's, thereby providing an extensive physical presence.\\u000b\\u000b3. Citigroup Inc.:\\u000bCitigroup is a globally diversified financial services holding company that provides a broad range of financial services to consumer and corporate customers. The financial giant is led by CEO Jane Fraser. Citigroup is more global in operations compared to JP Morgan, with business spanning over 100 countries. It also has a larger credit card business compared to JP Morgan.\"],\n  [13, \"**Ross McBride, Managing Director, Document and Business Solutions, JP Morgan Chase Battlecard:**\\u000b\\u000b1. **Persona Care-abouts:** Ross is likely most concerned with the efficiency and security of JP Morgan Chase's document and business solutions. As the leader in this space, he would prioritize maintaining the top-tier quality of service the firm and its clients expect. Staying on top of new innovations and trends in business solution strategies, particularly those that enhance operational efficiency and data security, would be crucial for him. He is also likely invested in supporting his team's effective performance and continuous professional development.\\u000b\\u000b2. **Challenges:** As part of a mega banking institution like JP Morgan Chase, Ross is faced with the continuous challenge of managing enormous volumes of data and resources. Given the increasing cyber-security threats, maintaining security for the firm's data systems is a key issue. Moreover, with the ongoing COVID-19 pandemic and the need for remote solutions, Ross would also need to deal with ensuring efficient, reliable, and secure remote working protocols and systems.\\u000b\\u000b3. **Ricoh's Alignment with their Needs:** Ricoh is a global leader in the provision of business services, and as such, offers a plethora of solutions to align with Ross's immediate needs and priorities. Our innovative document and data management solutions can assist in streamlining workflows, managing wide-ranging resources, and enhancing operational efficiency at JP Morgan Chase. We also have a strong thrust in security and are well equipped to provide robust, secure platforms reducing vulnerability to cyber threats. Plus, our expertise in workforce management and digital transformation makes us a strong partner in Ross's path to fostering team development and guiding the digital transformation journey at JP Morgan Chase.\\u000b\\u000bApproaching Ross should involve discussions on how the innovative solutions Ricoh provides can significantly boost his department's efficiency while preserving security and supporting his team's professional growth. Offering examples of proven track records in similar banking environments can further strengthen this strategy.\"],\n];\n\nfor (const [index, newText] of replacements) {\n  paragraphs.items[index].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the report-refresh edit: updated timestamp plus five rewritten\n# section bodies (Company Overview, Recent News, Competitors, Insights).\n# $vt (chr 11 / vertical tab) marks a Word manual line break (<w:br/>).\n# Each target paragraph's Range (minus its trailing paragraph mark) has\n# its .Text replaced wholesale, which keeps the paragraph (and its style)\n# in place while swapping the run content/line breaks inside it.\n$d = $word.ActiveDocument\n$vt = [char]11\n\n$replacements = @{\n  2 = 'Generated: 2025-08-06 14:02:06'\n  5 = 'As an AI developed by OpenAI, I can deliver information based on the data I was last trained on, which includes information up to September 2021. For this reason, I''m unable to provide real-time or precise company-specific information or updates that have occurred post-September 2021. However, I can provide a general response based on the details you have given above:' + $vt + '' + $vt + '1. JP Morgan Chase, one of the world''s most valuable banks, provides services to millions of consumers, small businesses, and many of the world''s most prominent corporate, institutional, and government clients. They''re major players in investment banking, financial services for consumers and businesses, financial transaction processing, asset management, and private equity.' + $vt + '' + $vt + '2. The current geopolitical environment, marked by the global economic recovery from the Covid-19 pandemic, trade tensions, and regulatory changes, can impact the bank''s business operations, risk profile, and market dynamics. These impacts can vary greatly depending on the specific circumstances.' + $vt + '' + $vt + '3. Amidst its competitors like Bank of America, Goldman Sachs, and Citigroup, JP Morgan Chase often stands out with its vast portfolio of services, robust financial performance, and a strong brand name. ' + $vt + '' + $vt + '4. Strengths of JP Morgan include its strong market position, diversified business model, and robust capital position. It also has a strong technology infrastructure that offers online and mobile banking options to customers.' + $vt + '' + $vt + '5. Key challenges for JP Morgan Chase can include regulatory changes, geopolitical uncertainties, and economic uncertainties related to global events like the Covid-19 pandemic. They undertake comprehensive risk management strategies to mitigate these risks, but the effectiveness of these strategies can vary depending on the scope of external events.' + $vt + '' + $vt + '6. JP Morgan Chase distinguishes itself from competitors with its strong emphasis on digital banking and technologies, a diversified business model, and strong risk management capabilities.' + $vt + '' + $vt + '7. I don''t have up-to-date data on recent projects or initiatives started by JP Morgan Chase. However, the company has been known for investing heavily in technological advancements and customer-centric initiatives.' + $vt + '' + $vt + '8. I am unable to provide specific forecast figures for future growth. It''s generally recommended to consult specialized financial sources or professionals for accurate figures.' + $vt + '' + $vt + '9. Comparison with competitors on revenue growth forecasts would also require up-to-date data from financial sources or professionals.' + $vt + '' + $vt + '10. Currently, I can''t provide highlights from the last investor day presentation as my last training data covers information up until September 2021.' + $vt + '' + $vt + '11. For detailed comparisons of the latest 10Qs, you should consult professional financial sources or regulatory disclosures. Even more so for the comparison with the latest 10Qs of other competitors.' + $vt + '' + $vt + '12. A digital transformation plan typically includes a broad range of initiatives that use digital technologies to modify or create new business processes, culture, and customer experiences to meet changing business and market requirements. As of my last training data, JP Morgan Chase has been actively investing in digital banking, but specific initiatives for the upcoming year are not known.' + $vt + '' + $vt + '13. Over the past years, JP Morgan Chase has implemented several digital transformation and operational efficiency initiatives. They have been investing heavily in technologies like AI and blockchain, improving their online and mobile banking platforms and customer-centric digital tools.' + $vt + '' + $vt + '14. The exact individuals responsible for digital transformation in JP Morgan Chase can change over time, but large corporations typically have a Chief Digital Officer and a team dedicated to digital transformation initiatives. They work in conjunction with other departments to ensure business-wide digital strategies align.' + $vt + '' + $vt + 'I encourage you to consult official JP Morgan Chase resources or related regulatory filings for more updated and accurate data.'\n  7 = '1. **Trump Claims Rejection by JP Morgan:** Donald Trump claimed that JP Morgan Chase and Bank of America rejected him as a customer, alleging discrimination (Source: CNBC, The New York Times, Politico, Financial Times). This could draw regulatory attention and potential legal implications for the firm, impacting its reputation and investor relations.' + $vt + '' + $vt + '2. **Employee Complaints about Fitness Center Costs:** JP Morgan staff have expressed dissatisfaction about having to pay for the new high-end fitness center at the company''s headquarters, as they return to office fulltime (Source: New York Post). This could possibly impact employee morale and turnover rates and might require the company to revisit its employee engagement strategies.' + $vt + '' + $vt + '3. **Expansion in Maryland:** JP Morgan plans to open four new branches in Maryland, signaling its intent to expand its footprint in the region through 2025 (Source: Maryland Daily Record). This is a positive development and could contribute to overall business growth and market share increase.' + $vt + '' + $vt + '4. **Acquisition of 3M Co. Shares:** The company has significantly expanded its stake in the 3M Co. with the purchase of over 13 million shares (Source: GuruFocus). This investment could enhance JP Morgan''s financial returns depending on the performance of 3M''s stock.' + $vt + '' + $vt + '5. **Potential Apple Play:** There are speculations that JP Morgan might make a play for Apple''s credit card portfolio, possibly signaling a greater ambition to conquer more business in technology-based financial services and products (Source: Payments Dive). This could symbolise strategic change, moving towards a more tech-oriented business model.' + $vt + ''\n  11 = '1. Bank of America Corp:' + $vt + 'Bank of America is one of the world''s largest financial institutions, serving individual consumers, small and middle-market businesses and large corporations with a broad spectrum of banking, investing, asset management and other financial services. Its finance and leadership segment is directed by CEO Brian Moynihan. One key differentiation between JP Morgan Chase and Bank of America is that the latter has a more diversified business mix and is heavily inclined towards consumer banking.' + $vt + '' + $vt + '2. Wells Fargo & Co:' + $vt + 'Wells Fargo is a diversified, community-based financial services company that provides banking, investment and mortgage products and services, as well as consumer and commercial finance. The company is helmed by CEO Charles Scharf. One key differentiator is its network of branches across the United States which is larger than JP Morgan''s, thereby providing an extensive physical presence.' + $vt + '' + $vt + '3. Citigroup Inc.:' + $vt + 'Citigroup is a globally diversified financial services holding company that provides a broad range of financial services to consumer and corporate customers. The financial giant is led by CEO Jane Fraser. Citigroup is more global in operations compared to JP Morgan, with business spanning over 100 countries. It also has a larger credit card business compared to JP Morgan.'\n  13 = '**Ross McBride, Managing Director, Document and Business Solutions, JP Morgan Chase Battlecard:**' + $vt + '' + $vt + '1. **Persona Care-abouts:** Ross is likely most concerned with the efficiency and security of JP Morgan Chase''s document and business solutions. As the leader in this space, he would prioritize maintaining the top-tier quality of service the firm and its clients expect. Staying on top of new innovations and trends in business solution strategies, particularly those that enhance operational efficiency and data security, would be crucial for him. He is also likely invested in supporting his team''s effective performance and continuous professional development.' + $vt + '' + $vt + '2. **Challenges:** As part of a mega banking institution like JP Morgan Chase, Ross is faced with the continuous challenge of managing enormous volumes of data and resources. Given the increasing cyber-security threats, maintaining security for the firm''s data systems is a key issue. Moreover, with the ongoing COVID-19 pandemic and the need for remote solutions, Ross would also need to deal with ensuring efficient, reliable, and secure remote working protocols and systems.' + $vt + '' + $vt + '3. **Ricoh''s Alignment with their Needs:** Ricoh is a global leader in the provision of business services, and as such, offers a plethora of solutions to align with Ross''s immediate needs and priorities. Our innovative document and data management solutions can assist in streamlining workflows, managing wide-ranging resources, and enhancing operational efficiency at JP Morgan Chase. We also have a strong thrust in security and are well equipped to provide robust, secure platforms reducing vulnerability to cyber threats. Plus, our expertise in workforce management and digital transformation makes us a strong partner in Ross''s path to fostering team development and guiding the digital transformation journey at JP Morgan Chase.' + $vt + '' + $vt + 'Approaching Ross should involve discussions on how the innovative solutions Ricoh provides can significantly boost his department''s efficiency while preserving security and supporting his team''s professional growth. Offering examples of proven track records in similar banking environments can further strengthen this strategy.'\n}\n\nforeach ($index in $replacements.Keys) {\n  # Word's Paragraphs collection is 1-based; $index is the 0-based\n  # paragraph number (matching the Office.js body.paragraphs.items index).\n  $para = $d.Paragraphs.Item([int]$index + 1)\n  $r = $para.Range\n  $r.MoveEnd(1, -1) | Out-Null\n  $r.Text = $replacements[$index]\n}\n"}
